$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 30: new participant sub_029, flagged with reason "unclear instructions"
$ws.Range("A30").Value = "sub_029"
$ws.Range("B30").Value = $true
$ws.Range("C30").Value = "unclear instructions"

# Row 31: new participant sub_030, not flagged
$ws.Range("A31").Value = "sub_030"
$ws.Range("B31").Value = $false

# Row 29: sub_028 is now flagged as a QC failure with a reason
$ws.Range("B29").Value = $true
$ws.Range("C29").Value = "extra participant"

$ws.Range("C29").Select()
